$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("C2").Value = 495872
$ws.Range("E2").Value = 764044962

# Row 19
$ws.Range("C19").Value = 117226
$ws.Range("E19").Value = 178790026

# Row 33
$ws.Range("C33").Value = 143111
$ws.Range("E33").Value = 222701540

# Row 45
$ws.Range("C45").Value = 103206
$ws.Range("E45").Value = 159030588

# Row 58
$ws.Range("C58").Value = 31011
$ws.Range("E58").Value = 50755569

# Row 67
$ws.Range("C67").Value = 216238
$ws.Range("E67").Value = 341218723

# Row 96
$ws.Range("C96").Value = 214593
$ws.Range("E96").Value = 323457894

# Row 111
$ws.Range("C111").Value = 857392
$ws.Range("E111").Value = 1399724697

# Row 152
$ws.Range("C152").Value = 132100
$ws.Range("D152").Value = 41347
$ws.Range("E152").Value = 206515908

# Row 164
$ws.Range("C164").Value = 350727
$ws.Range("E164").Value = 521029317

# Row 182
$ws.Range("C182").Value = 399651
$ws.Range("D182").Value = 122962
$ws.Range("E182").Value = 581422622

# Row 201
$ws.Range("E201").Value = 241519133

# Row 214
$ws.Range("C214").Value = 402275
$ws.Range("E214").Value = 605029400
